# RPA datasets push 2024-04-20
#
# The source IPO tracking sheet had two rows for "에이피알" (rows 10 and 16,
# 1-indexed incl. header) that were removed from the dataset. Deleting the
# rows shifts every row below each deletion up by one, collapsing the sheet
# from 17 data+header rows (A1:L17) down to 15 (A1:L15). Shared strings that
# become unused (the "2024-02-14" date, "에이피알" company name and the
# "신한, 하나" underwriter combo) are dropped automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the later row first so row 10's index isn't affected by the
# first deletion.
$ws.Rows(16).Delete()
$ws.Rows(10).Delete()
